$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1683
$ws1.Range("F5").Value = 368
$ws1.Range("F7").Value = 1113
$ws1.Range("F9").Value = 166
$ws1.Range("F10").Value = 166
$ws1.Range("F11").Value = 5
$ws1.Range("F12").Value = 1467
$ws1.Range("F13").Value = 3095
$ws1.Range("F14").Value = 639
$ws1.Range("F15").Value = 1782
$ws1.Range("F16").Value = 1809
$ws1.Range("F20").Value = 1482
$ws1.Range("F23").Value = 12
$ws1.Range("F24").Value = 1224
$ws1.Range("F25").Value = 410
$ws1.Range("F26").Value = 461
$ws1.Range("F27").Value = 126
$ws1.Range("F28").Value = 4820
$ws1.Range("F29").Value = 50
$ws1.Range("F30").Value = 755
$ws1.Range("F31").Value = 574
$ws1.Range("F32").Value = 1662
$ws1.Range("F34").Value = 151

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 58
$ws2.Range("F3").Value = 31
$ws2.Range("F5").Value = 29
$ws2.Range("F6").Value = 58
$ws2.Range("F7").Value = 78
$ws2.Range("G3").Value = 100

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 58
$ws4.Range("F5").Value = 31
$ws4.Range("F8").Value = 29
$ws4.Range("F9").Value = 58
$ws4.Range("F10").Value = 78
$ws4.Range("F12").Value = 1683
$ws4.Range("F13").Value = 368
$ws4.Range("F15").Value = 1113
$ws4.Range("F17").Value = 166
$ws4.Range("F18").Value = 166
$ws4.Range("F20").Value = 5
$ws4.Range("F21").Value = 1467
$ws4.Range("F22").Value = 3095
$ws4.Range("F23").Value = 639
$ws4.Range("F24").Value = 1782
$ws4.Range("F25").Value = 1809
$ws4.Range("F29").Value = 1482
$ws4.Range("F33").Value = 12
$ws4.Range("F35").Value = 1224
$ws4.Range("F36").Value = 410
$ws4.Range("F37").Value = 461
$ws4.Range("F38").Value = 126
$ws4.Range("F39").Value = 4820
$ws4.Range("F40").Value = 50
$ws4.Range("F41").Value = 755
$ws4.Range("F42").Value = 574
$ws4.Range("F43").Value = 1662
$ws4.Range("F47").Value = 151
$ws4.Range("G5").Value = 100
